$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values (Price / Volume(1h) columns, and
# the KickToken/CEJI row swap in columns B/C) while preserving the
# original literal-text representation (avoid Excel auto-converting
# numeric-looking / percent-looking strings into real numbers).
$cellValues = @{
    'D2' = '245.39'
    'E2' = '-0.50%'
    'D3' = '29.82'
    'E3' = '0.01%'
    'D4' = '5.211'
    'E4' = '0.87%'
    'D5' = '0.05747'
    'E5' = '0.61%'
    'D6' = '6.645'
    'E6' = '0.52%'
    'D7' = '3.275'
    'E7' = '7.34%'
    'D8' = '0.8581'
    'E8' = '0.01%'
    'D9' = '0.8526'
    'E9' = '-1.67%'
    'D10' = '0.1382'
    'E10' = '1.06%'
    'D11' = '0.07087'
    'E11' = '0.24%'
    'D12' = '0.03145'
    'E12' = '9.91%'
    'D13' = '0.09352'
    'E13' = '-0.40%'
    'D14' = '0.001523'
    'E14' = '-0.69%'
    'D15' = '0.0005948'
    'E15' = '-94.22%'
    'D16' = '0.005965'
    'E16' = '-2.83%'
    'E17' = '1.33%'
    'D18' = '2.193'
    'E18' = '-3.33%'
    'E19' = '0.36%'
    'D20' = '0.03332'
    'E20' = '2.13%'
    'E21' = '0.58%'
    'D22' = '3.479'
    'E22' = '19.69%'
    'D23' = '0.1409'
    'E23' = '2.08%'
    'D24' = '0.04115'
    'E24' = '-1.01%'
    'E25' = '0.19%'
    'E26' = '-18.19%'
    'D28' = '0.0001448'
    'E28' = '-25.29%'
    'D40' = '0.03762'
    'E40' = '-0.15%'
    'D41' = '0.1070'
    'E41' = '-0.01%'
    'B42' = 'CEJI'
    'C42' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'D42' = '0.002459'
    'E42' = '17.08%'
    'B43' = 'KickToken'
    'C43' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'D43' = '0.002949'
    'E43' = '-48.37%'
    'D44' = '0.01025'
    'E44' = '4.51%'
    'D45' = '0.00005465'
    'E45' = '6.99%'
    'D47' = '0.08985'
    'E47' = '19.62%'
    'D48' = '0.002209'
    'E48' = '-19.03%'
    'D49' = '0.00002099'
    'D50' = '0.0001999'
}

foreach ($addr in $cellValues.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cellValues[$addr]
    $rng.Style = "Normal"
}
